$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.968.95"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").Value = "1.983.01"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.64"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("E6").Value = "  +1.87%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.17"
$ws.Range("E7").Value = "  +3.34%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.384"
$ws.Range("E9").Value = "  +2.85%  "
$ws.Range("E10").Value = "  -1.51%  "
$ws.Range("E11").Value = "  +0.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.80"
$ws.Range("E12").Value = "  +7.80%  "
$ws.Range("E13").Value = "  +2.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.21"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").Value = "2.271.51"
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("E16").Value = "  +3.91%  "
$ws.Range("D17").Value = "1.975.91"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").Value = "36.843.02"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.24"
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("D20").Value = "0.0₃0865"
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("E21").Value = "  +1.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.40"
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.51"
$ws.Range("E24").Value = "  +2.67%  "
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("E26").Value = "  +1.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.33"
$ws.Range("E27").Value = "  +0.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.92"
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.57"
$ws.Range("E29").Value = "  +0.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.35"
$ws.Range("E30").Value = "  +18.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.122"
$ws.Range("E31").Value = "  +1.93%  "
$ws.Range("E32").Value = "  +2.93%  "
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.53"
$ws.Range("E34").Value = "  +5.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.28"
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.51"
$ws.Range("E39").Value = "  -8.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.100"
$ws.Range("E40").Value = "  +1.54%  "
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("E42").Value = "  +0.30%  "
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.38"
$ws.Range("E44").Value = "  +1.11%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.19"
$ws.Range("E45").Value = "  +2.67%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.370.41"
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.27"
$ws.Range("E48").Value = "  +1.47%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.00"
$ws.Range("E49").Value = "  +12.94%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "46.34"
$ws.Range("E50").Value = "  +5.36%  "
$ws.Range("B51").Value = "MXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.81"
$ws.Range("E51").Value = "  -0.75%  "
